$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-25 Sunday" "2024-02-26 Monday"

Replace-Text "202×2=" "400×5="
Replace-Text "930×5=" "882×2="
Replace-Text "448×3=" "452×7="
Replace-Text "923×6=" "311×9="
Replace-Text "449×9=" "355×4="

Replace-Text "509×8=" "256×2="
Replace-Text "704×5=" "932×6="
Replace-Text "874×7=" "325×7="
Replace-Text "171×9=" "792×7="
Replace-Text "819×7=" "190×6="

Replace-Text "151×8=" "452×5="
Replace-Text "937×2=" "532×8="
Replace-Text "512×3=" "196×3="
Replace-Text "790×8=" "157×6="
Replace-Text "622×3=" "533×5="

Replace-Text "406×5=" "218×9="
Replace-Text "177×4=" "799×2="
Replace-Text "224×7=" "982×3="
Replace-Text "123×7=" "274×2="
Replace-Text "337×5=" "717×6="

Replace-Text "537×7=" "119×6="
Replace-Text "236×3=" "223×5="
Replace-Text "264×2=" "633×9="
Replace-Text "999×4=" "264×4="
Replace-Text "667×9=" "477×4="
